$wb = $excel.ActiveWorkbook

# --- Receitas sheet: fill in new revenue rows (client invoices for "GRUPO XYZ") ---
$receitas = $wb.Worksheets.Item("Receitas")

$receitas.Range("A2").Value = 45902
$receitas.Range("B2").Value = 150000
$receitas.Range("C2").Value = "4567-8"
$receitas.Range("D2").Value = "GRUPO XYZ"
$receitas.Range("E2").Value = "CLIENTE DO GRUPO XYZ LTDA"
$receitas.Range("I2").Value = 1

$receitas.Range("A3").Value = 45930
$receitas.Range("B3").Value = 50000
$receitas.Range("C3").Value = "4567-8"
$receitas.Range("D3").Value = "GRUPO XYZ"
$receitas.Range("E3").Value = "CLIENTE DO GRUPO XYZ LTDA"
$receitas.Range("I3").Value = 2

$receitas.Range("A4").Value = 45901
$receitas.Range("B4").Value = 350000
$receitas.Range("C4").Value = "8765-4"
$receitas.Range("D4").Value = "GRUPO XYZ"
$receitas.Range("E4").Value = "CLIENTE DO GRUPO XYZ LTDA"
$receitas.Range("I4").Value = 3

$receitas.Range("A5").Value = 45910
$receitas.Range("B5").Value = 50000
$receitas.Range("C5").Value = "8765-4"
$receitas.Range("D5").Value = "GRUPO XYZ"
$receitas.Range("E5").Value = "CLIENTE DO GRUPO XYZ LTDA"
$receitas.Range("I5").Value = 4

# --- Tarifas bancárias: leave selection parked on D7 ---
$tarifas = $wb.Worksheets.Item("Tarifas bancárias")
$tarifas.Activate()
$tarifas.Range("D7").Select()

# --- Make "Receitas" the active/selected tab, cursor on I1 ---
$receitas.Activate()
$receitas.Range("I1").Select()
